$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 203.7816646666667
$ws.Range("H2").Value = 611.344994
$ws.Range("I2").Value = 0.6667327591988204
$ws.Range("J2").Value = 0.6667327591988205
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7651113333333334
$ws.Range("N2").Value = 2.295334
$ws.Range("O2").Value = 0.1330257938600752
$ws.Range("P2").Value = 0.1330257938600752
$ws.Range("Q2").Value = 155.9156611619996
$ws.Range("R2").Value = 1403.240950457996
$ws.Range("S2").Value = 0.0886926545849414
$ws.Range("T2").Value = 0.08869265458494142

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 203.7816646666667
$ws.Range("H3").Value = 611.344994
$ws.Range("I3").Value = 0.6667327591988204
$ws.Range("J3").Value = 0.6667327591988205
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8520123333333333
$ws.Range("N3").Value = 2.556037
$ws.Range("O3").Value = 0.1481348035016799
$ws.Range("P3").Value = 0.1481348035016799
$ws.Range("Q3").Value = 173.6244916031976
$ws.Range("R3").Value = 1562.620424428778
$ws.Range("S3").Value = 0.09876632627205009
$ws.Range("T3").Value = 0.09876632627205012

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 203.7816646666667
$ws.Range("H4").Value = 611.344994
$ws.Range("I4").Value = 0.6667327591988204
$ws.Range("J4").Value = 0.6667327591988205
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.134477666666666
$ws.Range("N4").Value = 12.403433
$ws.Range("O4").Value = 0.7188394026382449
$ws.Range("P4").Value = 0.718839402638245
$ws.Range("Q4").Value = 842.530741440489
$ws.Range("R4").Value = 7582.776672964403
$ws.Range("S4").Value = 0.4792737783418288
$ws.Range("T4").Value = 0.4792737783418289

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.2065833519051582
$ws.Range("J5").Value = 0.2065833519051582
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7651113333333334
$ws.Range("N5").Value = 2.295334
$ws.Range("O5").Value = 0.1330257938600752
$ws.Range("P5").Value = 0.1330257938600752
$ws.Range("Q5").Value = 48.30958049227911
$ws.Range("R5").Value = 434.786224430512
$ws.Range("S5").Value = 0.02748091438545893
$ws.Range("T5").Value = 0.02748091438545894

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("H6").Value = 189.421768
$ws.Range("I6").Value = 0.2065833519051582
$ws.Range("J6").Value = 0.2065833519051582
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8520123333333333
$ws.Range("N6").Value = 2.556037
$ws.Range("O6").Value = 0.1481348035016799
$ws.Range("P6").Value = 0.1481348035016799
$ws.Range("Q6").Value = 53.7965608459351
$ws.Range("R6").Value = 484.169047613416
$ws.Range("S6").Value = 0.03060218424118899
$ws.Range("T6").Value = 0.03060218424118899

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("H7").Value = 189.421768
$ws.Range("I7").Value = 0.2065833519051582
$ws.Range("J7").Value = 0.2065833519051582
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.134477666666666
$ws.Range("N7").Value = 12.403433
$ws.Range("O7").Value = 0.7188394026382449
$ws.Range("P7").Value = 0.718839402638245
$ws.Range("Q7").Value = 261.0533564588382
$ws.Range("R7").Value = 2349.480208129544
$ws.Range("S7").Value = 0.1485002532785102
$ws.Range("T7").Value = 0.1485002532785103

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 38.719942
$ws.Range("H8").Value = 116.159826
$ws.Range("I8").Value = 0.1266838888960214
$ws.Range("J8").Value = 0.1266838888960214
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.7651113333333334
$ws.Range("N8").Value = 2.295334
$ws.Range("O8").Value = 0.1330257938600752
$ws.Range("P8").Value = 0.1330257938600752
$ws.Range("Q8").Value = 29.62506645020933
$ws.Range("R8").Value = 266.625598051884
$ws.Range("S8").Value = 0.01685222488967481
$ws.Range("T8").Value = 0.01685222488967481

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 38.719942
$ws.Range("H9").Value = 116.159826
$ws.Range("I9").Value = 0.1266838888960214
$ws.Range("J9").Value = 0.1266838888960214
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8520123333333333
$ws.Range("N9").Value = 2.556037
$ws.Range("O9").Value = 0.1481348035016799
$ws.Range("P9").Value = 0.1481348035016799
$ws.Range("Q9").Value = 32.98986812995133
$ws.Range("R9").Value = 296.908813169562
$ws.Range("S9").Value = 0.01876629298844077
$ws.Range("T9").Value = 0.01876629298844078

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 38.719942
$ws.Range("H10").Value = 116.159826
$ws.Range("I10").Value = 0.1266838888960214
$ws.Range("J10").Value = 0.1266838888960214
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.134477666666666
$ws.Range("N10").Value = 12.403433
$ws.Range("O10").Value = 0.7188394026382449
$ws.Range("P10").Value = 0.718839402638245
$ws.Range("Q10").Value = 160.0867354536286
$ws.Range("R10").Value = 1440.780619082658
$ws.Range("S10").Value = 0.0910653710179058
$ws.Range("T10").Value = 0.09106537101790584

Write-Output "applied updates"
